$wb = $excel.ActiveWorkbook

# --- Sheet 1: Batter ---
$ws = $wb.Worksheets.Item("Batter")

$ws.Range("C2").Value = 1.694062797813139
$ws.Range("E2").Value = 0.742319372979656

$ws.Range("C3").Value = 2.225099032184885
$ws.Range("E3").Value = 0.5554491472675398

$ws.Range("C4").Value = 1.940739869809948
$ws.Range("E4").Value = 0.6618126070274846

$ws.Range("C5").Value = 2.029920287197169
$ws.Range("D5").Value = 0.9884782555672129
$ws.Range("E5").Value = 0.6300178890845086

$ws.Range("C6").Value = 2.171775041616835
$ws.Range("E6").Value = 0.5765009566294313

$ws.Range("C7").Value = 2.178618141920127
$ws.Range("D7").Value = 0.9170390054603687
$ws.Range("E7").Value = 0.5738279245335569

$ws.Range("C8").Value = 2.0469229793699
$ws.Range("D8").Value = 0.9821102070623208
$ws.Range("E8").Value = 0.6237939624380013

# --- Sheet 2: Pitcher ---
$ws2 = $wb.Worksheets.Item("Pitcher")

$ws2.Range("C2").Value = 15.01932178290843
$ws2.Range("E2").Value = -1.01570273264918

$ws2.Range("C3").Value = 6.683836766917238
$ws2.Range("E3").Value = 0.6008119695625709

$ws2.Range("C4").Value = 8.169039967069054
$ws2.Range("E4").Value = 0.4036957806847478

$ws2.Range("C5").Value = 5.509507065597373
$ws2.Range("D5").Value = 0.9911188452018577
$ws2.Range("E5").Value = 0.7287615029017056

$ws2.Range("C6").Value = 6.505615417033987
$ws2.Range("E6").Value = 0.6218164685608235

$ws2.Range("C7").Value = 6.202958542276117
$ws2.Range("D7").Value = 0.9050747766293435
$ws2.Range("E7").Value = 0.6561859678006043

$ws2.Range("C8").Value = 5.906910872072573
$ws2.Range("D8").Value = 0.9890687624604397
$ws2.Range("E8").Value = 0.6882211315531173

$wb.Save()
